# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q4" (per-fund holdings) and
# "总计" (totals-by-quarter). This adds a new "2022-Q1" per-fund holdings
# sheet (cloned from "2021-Q4" so it keeps identical layout/formatting),
# placed between the two existing sheets, and records the new quarter at
# the top of the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" sheet by duplicating "2021-Q4" (keeps the
#    header row + styling identical) and placing it right after it.
#    NOTE: worksheet handles captured before a Copy() can end up bound
#    to the wrong tab once the sheet collection shifts, so re-fetch
#    "总计" by name again afterwards instead of reusing an old handle.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Columns B (fund code) and D:G (size/position numbers stored as text in
# the source data) must stay text so values like "002236" keep their
# leading zero instead of becoming the number 2236. Mark them as Text
# before writing, then clear the format again so the cells end up with
# no explicit style (matching the rest of the sheet) while keeping the
# text type.
$q1.Range("B2:B3").NumberFormat = "@"
$q1.Range("D2:G3").NumberFormat = "@"

$q1.Range("B2").Value = "002236"
$q1.Range("C2").Value = "大成中证360互联网+大数据100指数A"
$q1.Range("D2").Value = "5.67"
$q1.Range("E2").Value = "93.32"
$q1.Range("F2").Value = "1.09"
$q1.Range("G2").Value = "0.0618"
$q1.Range("H2").Value = 2

$q1.Range("B3").Value = "003359"
$q1.Range("C3").Value = "大成中证360互联网+大数据100指数C"
$q1.Range("D3").Value = "4.08"
$q1.Range("E3").Value = "93.32"
$q1.Range("F3").Value = "1.09"
$q1.Range("G3").Value = "0.0445"
$q1.Range("H3").Value = 2

$q1.Range("B2:B3").ClearFormats()
$q1.Range("D2:G3").ClearFormats()

# ---------------------------------------------------------------------
# 2) Push the existing "总计" row down and insert the 2022-Q1 totals
#    above it, so the newest quarter is listed first.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

$total.Range("A3").Value = 1
